$wb = $excel.ActiveWorkbook

# Duplicate the last sheet (Bus_Makulu_r) - it carries the same layout,
# styles, dxfs and conditional formatting used by the new droplink sheet.
$srcSheet = $wb.Worksheets.Item("Bus_Makulu_r")
$srcSheet.Copy([System.Reflection.Missing]::Value, $srcSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Trailer1Axle_f"

# Match the view/pane setup used by the other "_f" sheets.
$newSheet.Range("H7").Select() | Out-Null

# Row label (H3) identifies this droplink instance.
$newSheet.Range("H3").Value = "Droplink_Trailer1Axle_f"

# Updated hardpoints / values for the new Trailer1Axle_f droplink.
$newSheet.Range("F5").Value = 0.05
$newSheet.Range("G5").Value = 0.6
$newSheet.Range("H5").Value = 0.19

$newSheet.Range("F6").Formula = "=0.3-0.15"
$newSheet.Range("G6").Value = 0.57999999999999996
$newSheet.Range("H6").Value = 0.2

$newSheet.Range("H7").Value = 50
$newSheet.Range("H8").Value = 0.5
